$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 406.7620641209767
$ws.Range("C2").Value = 1423.784512505422
$ws.Range("D2").Value = 1826.697150462707
$ws.Range("E2").Value = 64.62950279638585
$ws.Range("F2").Value = 19.42673116764596
$ws.Range("G2").Value = 307.1297081676828
$ws.Range("H2").Value = 26.3382283904123
$ws.Range("I2").Value = 688.7475634103021
$ws.Range("J2").Value = 16.35272369515465
$ws.Range("K2").Value = 5.513404951952301
$ws.Range("L2").Value = 117.8835876429822
$ws.Range("M2").Value = 1.221204087471981
$ws.Range("N2").Value = 1.006854132101349
$ws.Range("O2").Value = 13.39270749288584
$ws.Range("Q2").Value = 64.72721744196113
$ws.Range("S2").Value = 64.72721744196113
$ws.Range("T2").Value = 64.72721744196113
$ws.Range("U2").Value = 1.633954832474227
$ws.Range("Z2").Value = 1131.535015516773

$ws.Range("B3").Value = 17861.818533069
$ws.Range("C3").Value = 124276.8874610486
$ws.Range("D3").Value = 95598.60442291811
$ws.Range("E3").Value = 4533.759471289621
$ws.Range("F3").Value = 5762.283039516643
$ws.Range("G3").Value = 6432.12558334999
$ws.Range("H3").Value = 17.00465880233554
$ws.Range("I3").Value = 168683.7684859641
$ws.Range("J3").Value = 292.2271643795879
$ws.Range("K3").Value = 4.255912753998186
$ws.Range("L3").Value = 2472.298894748137
$ws.Range("M3").Value = 64.04697408836569
$ws.Range("N3").Value = 267.0271231858534
$ws.Range("O3").Value = 151.3571743724798
$ws.Range("Q3").Value = 11259.30068662367
$ws.Range("S3").Value = 11259.30068662367
$ws.Range("T3").Value = 11259.30068662367
$ws.Range("U3").Value = 1.182775708762887
$ws.Range("Z3").Value = 75058.96728275302

$ws.Range("B4").Value = 417.0983106078516
$ws.Range("C4").Value = 730.3517250623604
$ws.Range("D4").Value = 1673.703985177697
$ws.Range("E4").Value = 70.05364113557339
$ws.Range("F4").Value = 12.32842606886002
$ws.Range("G4").Value = 321.2683864559366
$ws.Range("H4").Value = 23.18353026380176
$ws.Range("I4").Value = 463.6137425946026
$ws.Range("J4").Value = 16.60971840232479
$ws.Range("K4").Value = 5.997600285112068
$ws.Range("L4").Value = 116.8732224632952
$ws.Range("M4").Value = 1.250543040696576
$ws.Range("N4").Value = 0.7988417300089758
$ws.Range("O4").Value = 15.58366849701403
$ws.Range("Q4").Value = 53.17998944211951
$ws.Range("S4").Value = 53.17998944211951
$ws.Range("T4").Value = 53.17998944211951
$ws.Range("U4").Value = 1.903133956185567
$ws.Range("Z4").Value = 1095.504932171302

$ws.Range("B5").Value = 409.595483460815
$ws.Range("C5").Value = 1291.95929768566
$ws.Range("D5").Value = 1740.755188175594
$ws.Range("E5").Value = 70.04260997750842
$ws.Range("F5").Value = 36.88105418043184
$ws.Range("G5").Value = 284.1467258585485
$ws.Range("H5").Value = 19.26651188320137
$ws.Range("I5").Value = 1277.309704408859
$ws.Range("J5").Value = 14.39328545874083
$ws.Range("K5").Value = 3.979739452712001
$ws.Range("L5").Value = 92.2904584315487
$ws.Range("M5").Value = 1.16362985944985
$ws.Range("N5").Value = 1.769198074903962
$ws.Range("O5").Value = 11.200398795048
$ws.Range("Q5").Value = 109.7489180334657
$ws.Range("S5").Value = 109.7489180334657
$ws.Range("T5").Value = 109.7489180334657
$ws.Range("U5").Value = 1.826277963917526
$ws.Range("Z5").Value = 1061.204673530374
